$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add pooled scalar values to column E (ScalarPooled)
$ws.Range("E2").Value = 5.19
$ws.Range("E3").Value = 3.42
$ws.Range("E4").Value = 2.38
$ws.Range("E5").Value = 11.78

# Update workbook window position
$wb.Windows.Item(1).Left = 4800
$wb.Windows.Item(1).Top = 2280
